# Insert a new (empty) paragraph right after the paragraph that ends with
# "(/usr/themes && /usr/icons)." — i.e. immediately before the blank
# paragraph that already follows it (so the blank paragraph's own
# formatting - left indent 220, Times New Roman 18pt, justified - is what
# the newly inserted paragraph inherits, matching the target XML exactly).

$d = $word.ActiveDocument

# Locate the target sentence.
$search = $d.Content
$found = $search.Find.Execute("(/usr/themes && /usr/icons).", $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Paragraph containing the found text.
    $targetPara = $search.Paragraphs(1)

    # The paragraph that immediately follows it (currently the first of the
    # two pre-existing blank paragraphs).
    $followingPara = $targetPara.Next()

    # Collapse to the very start of that following paragraph and insert a
    # brand-new blank paragraph right before it, so it inherits the
    # following paragraph's formatting (left indent 220 etc.).
    $insertionPoint = $followingPara.Range
    $insertionPoint.Collapse(1)
    $insertionPoint.InsertParagraphBefore()
}
